# Insert a new data row at row 305 (pushing existing rows 305-423 down to 306-424)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 305; Excel shifts rows 305:423 -> 306:424
# and keeps formatting consistent with the surrounding rows (mirrors the
# canonical diff, which shows every row from 305 onward shifting down by one
# and a brand-new row 305 appearing with fresh data).
$ws.Rows(305).Insert()

# Populate the newly inserted row 305 with its data.
$ws.Range("A305").Value = 3
$ws.Range("B305").Value = "Femacal de La Calera"
$ws.Range("C305").Value = "Coquimbo"
$ws.Range("D305").Value = 45009
$ws.Range("E305").Value = 5
$ws.Range("F305").Value = 100112001
$ws.Range("G305").Value = "Berenjena"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 105
$ws.Range("K305").Value = 7000
$ws.Range("L305").Value = 7500
$ws.Range("M305").Value = 7262
$ws.Range("N305").Value = "$/caja 60 unidades"
$ws.Range("O305").Value = "Región de Arica y Parinacota"
$ws.Range("P305").Value = 121
$ws.Range("Q305").Value = 60
$ws.Range("R305").Value = "Hortaliza"
